# Revised gameplay scripts - mark resolved / superseded usage notes on the
# "Scripts" sheet (several Script entries got replaced by toolkit
# interactables during the interaction pass, so their old notes are
# annotated as deleted/no-longer-relevant), narrow the notes column now
# that it holds shorter text, and leave the selection where the author
# was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scripts")

$ws.Range("D14").Value = $ws.Range("D14").Text + " (deleted)"
$ws.Range("C15").Value = $ws.Range("C15").Text + " (NR)"
$ws.Range("D16").Value = $ws.Range("D16").Text + " (deleted)"
$ws.Range("D17").Value = $ws.Range("D17").Text + " (deleted)"
$ws.Range("D18").Value = $ws.Range("D18").Text + " (deleted"
$ws.Range("D19").Value = $ws.Range("D19").Text + " (deleted)"

# Narrow the "Usage Notes" column (was 56.43, now ~42.57 characters wide).
$ws.Columns.Item(3).ColumnWidth = 41.65

# Restore the cursor to the cell the author last edited.
$ws.Range("D14").Select()
